$d = $word.ActiveDocument

# Find the paragraph that ends with "Flora Class added" - the new bullet
# point is inserted immediately after it.
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Flora Class added") {
        $targetIndex = $i
    }
    $i = $i + 1
}

if ($targetIndex -eq -1) {
    $targetIndex = $d.Paragraphs.Count
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a new paragraph right after it; it inherits the same
# list-paragraph formatting (ListParagraph style, bullet numbering,
# Times New Roman 12pt) from the preceding paragraph automatically.
$target.Range.InsertParagraphAfter()

# Re-query the collection (rather than holding a stale reference) to get
# the freshly inserted paragraph and fill in its text.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "0.1% for bush to grow"
